$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the Price (column D) and Volume(1h) (column E) values for
# each coin row with the latest scrape from the GitHub Actions job.

$ws.Cells.Item(2, 4).Value = "23.940.18"
$ws.Cells.Item(2, 5).Value = "  +0.14%  "
$ws.Cells.Item(3, 4).Value = "1.651.43"
$ws.Cells.Item(3, 5).Value = "  +0.25%  "
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "0.9995"
$ws.Cells.Item(4, 5).Value = "  -0.69%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "311.04"
$ws.Cells.Item(5, 5).Value = "  +0.59%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.9992"
$ws.Cells.Item(6, 5).Value = "  -0.69%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.3896"
$ws.Cells.Item(7, 5).Value = "  -0.54%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.3844"
$ws.Cells.Item(8, 5).Value = "  -0.33%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "51.07"
$ws.Cells.Item(9, 5).Value = "  -0.34%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "1.347"
$ws.Cells.Item(10, 5).Value = "  -1.00%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.9994"
$ws.Cells.Item(11, 5).Value = "  -0.67%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.08449"
$ws.Cells.Item(12, 5).Value = "  -0.31%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "23.95"
$ws.Cells.Item(13, 5).Value = "  +0.07%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "7.031"
$ws.Cells.Item(14, 5).Value = "  -2.45%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "7.900"
$ws.Cells.Item(15, 5).Value = "  +0.34%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "0.00001318"
$ws.Cells.Item(16, 5).Value = "  +0.52%  "
$ws.Cells.Item(17, 4).Value = "1.653.15"
$ws.Cells.Item(17, 5).Value = "  +0.43%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "94.03"
$ws.Cells.Item(18, 5).Value = "  -0.70%  "
$ws.Cells.Item(19, 5).Value = "  -0.38%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "19.59"
$ws.Cells.Item(20, 5).Value = "  -2.27%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "6.941"
$ws.Cells.Item(21, 5).Value = "  +0.02%  "
$ws.Cells.Item(22, 5).Value = "  -0.64%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "13.65"
$ws.Cells.Item(23, 5).Value = "  +0.27%  "
$ws.Cells.Item(24, 4).Value = "23.935.33"
$ws.Cells.Item(24, 5).Value = "  +0.09%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "2.437"
$ws.Cells.Item(25, 5).Value = "  -0.40%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "2.920"
$ws.Cells.Item(26, 5).Value = "  -5.53%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "21.99"
$ws.Cells.Item(27, 5).Value = "  -0.81%  "
$ws.Cells.Item(28, 5).Value = "  -0.68%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "5.391"
$ws.Cells.Item(29, 5).Value = "  +1.34%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "137.37"
$ws.Cells.Item(30, 5).Value = "  -1.52%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "7.735"
$ws.Cells.Item(31, 5).Value = "  -1.48%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "2.486"
$ws.Cells.Item(32, 5).Value = "  -0.48%  "
$ws.Cells.Item(33, 4).Value = "1.834.66"
$ws.Cells.Item(33, 5).Value = "  +0.19%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "0.08113"
$ws.Cells.Item(34, 5).Value = "  +0.02%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "0.9926"
$ws.Cells.Item(35, 5).Value = "  -3.71%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.02934"
$ws.Cells.Item(36, 5).Value = "  -2.79%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "6.720"
$ws.Cells.Item(37, 5).Value = "  +0.58%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.2693"
$ws.Cells.Item(38, 5).Value = "  -0.51%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "10.50"
$ws.Cells.Item(39, 5).Value = "  -4.32%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.09120"
$ws.Cells.Item(40, 5).Value = "  -0.36%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.7565"
$ws.Cells.Item(41, 5).Value = "  +0.45%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "13.44"
$ws.Cells.Item(42, 5).Value = "  -0.21%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "1.423"
$ws.Cells.Item(43, 5).Value = "  -0.05%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "16.68"
$ws.Cells.Item(44, 5).Value = "  +2.55%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.6940"
$ws.Cells.Item(45, 5).Value = "  +0.09%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "2.446"
$ws.Cells.Item(46, 5).Value = "  -1.38%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "4.101"
$ws.Cells.Item(47, 5).Value = "  +0.29%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "1.001"
$ws.Cells.Item(48, 5).Value = "  -0.51%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.08287"
$ws.Cells.Item(49, 5).Value = "  +0.15%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "134.34"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "1.227"
$ws.Cells.Item(51, 5).Value = "  -0.16%  "
